$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '37.902.05'
$ws.Cells.Item(2, 4).Style = 'Normal'
$ws.Cells.Item(2, 5).NumberFormat = '@'
$ws.Cells.Item(2, 5).Value = '  +6.63%  '
$ws.Cells.Item(2, 5).Style = 'Normal'
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '2.058.74'
$ws.Cells.Item(3, 4).Style = 'Normal'
$ws.Cells.Item(3, 5).NumberFormat = '@'
$ws.Cells.Item(3, 5).Value = '  +3.89%  '
$ws.Cells.Item(3, 5).Style = 'Normal'
$ws.Cells.Item(4, 5).NumberFormat = '@'
$ws.Cells.Item(4, 5).Value = '  -0.11%  '
$ws.Cells.Item(4, 5).Style = 'Normal'
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '253.33'
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).NumberFormat = '@'
$ws.Cells.Item(5, 5).Value = '  +4.54%  '
$ws.Cells.Item(5, 5).Style = 'Normal'
$ws.Cells.Item(6, 5).NumberFormat = '@'
$ws.Cells.Item(6, 5).Value = '  +2.79%  '
$ws.Cells.Item(6, 5).Style = 'Normal'
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '65.60'
$ws.Cells.Item(7, 4).Style = 'Normal'
$ws.Cells.Item(7, 5).NumberFormat = '@'
$ws.Cells.Item(7, 5).Value = '  +14.96%  '
$ws.Cells.Item(7, 5).Style = 'Normal'
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.999'
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Cells.Item(8, 5).NumberFormat = '@'
$ws.Cells.Item(8, 5).Value = '  -0.14%  '
$ws.Cells.Item(8, 5).Style = 'Normal'
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '60.89'
$ws.Cells.Item(9, 4).Style = 'Normal'
$ws.Cells.Item(9, 5).NumberFormat = '@'
$ws.Cells.Item(9, 5).Value = '  +2.32%  '
$ws.Cells.Item(9, 5).Style = 'Normal'
$ws.Cells.Item(10, 5).NumberFormat = '@'
$ws.Cells.Item(10, 5).Value = '  +5.71%  '
$ws.Cells.Item(10, 5).Style = 'Normal'
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.0764'
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).NumberFormat = '@'
$ws.Cells.Item(11, 5).Value = '  +4.85%  '
$ws.Cells.Item(11, 5).Style = 'Normal'
$ws.Cells.Item(12, 5).NumberFormat = '@'
$ws.Cells.Item(12, 5).Value = '  +2.19%  '
$ws.Cells.Item(12, 5).Style = 'Normal'
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '0.921'
$ws.Cells.Item(13, 4).Style = 'Normal'
$ws.Cells.Item(13, 5).NumberFormat = '@'
$ws.Cells.Item(13, 5).Value = '  +0.33%  '
$ws.Cells.Item(13, 5).Style = 'Normal'
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '14.99'
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(14, 5).NumberFormat = '@'
$ws.Cells.Item(14, 5).Value = '  +6.88%  '
$ws.Cells.Item(14, 5).Style = 'Normal'
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '2.358.56'
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).NumberFormat = '@'
$ws.Cells.Item(15, 5).Value = '  +3.70%  '
$ws.Cells.Item(15, 5).Style = 'Normal'
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '20.75'
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Cells.Item(16, 5).NumberFormat = '@'
$ws.Cells.Item(16, 5).Value = '  +20.85%  '
$ws.Cells.Item(16, 5).Style = 'Normal'
$ws.Cells.Item(17, 5).NumberFormat = '@'
$ws.Cells.Item(17, 5).Value = '  +6.69%  '
$ws.Cells.Item(17, 5).Style = 'Normal'
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '2.047.76'
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Cells.Item(18, 5).NumberFormat = '@'
$ws.Cells.Item(18, 5).Value = '  +2.88%  '
$ws.Cells.Item(18, 5).Style = 'Normal'
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '37.770.63'
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).NumberFormat = '@'
$ws.Cells.Item(19, 5).Value = '  +6.40%  '
$ws.Cells.Item(19, 5).Style = 'Normal'
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '74.06'
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Cells.Item(20, 5).NumberFormat = '@'
$ws.Cells.Item(20, 5).Value = '  +4.76%  '
$ws.Cells.Item(20, 5).Style = 'Normal'
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '0.0₃0880'
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).NumberFormat = '@'
$ws.Cells.Item(21, 5).Value = '  +5.29%  '
$ws.Cells.Item(21, 5).Style = 'Normal'
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '5.38'
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).NumberFormat = '@'
$ws.Cells.Item(22, 5).Value = '  +6.37%  '
$ws.Cells.Item(22, 5).Style = 'Normal'
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '239.62'
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).NumberFormat = '@'
$ws.Cells.Item(23, 5).Value = '  +2.70%  '
$ws.Cells.Item(23, 5).Style = 'Normal'
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '2.69'
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Cells.Item(24, 5).NumberFormat = '@'
$ws.Cells.Item(24, 5).Value = '  +15.55%  '
$ws.Cells.Item(24, 5).Style = 'Normal'
$ws.Cells.Item(25, 5).NumberFormat = '@'
$ws.Cells.Item(25, 5).Value = '  -0.04%  '
$ws.Cells.Item(25, 5).Style = 'Normal'
$ws.Cells.Item(26, 5).NumberFormat = '@'
$ws.Cells.Item(26, 5).Value = '  +5.25%  '
$ws.Cells.Item(26, 5).Style = 'Normal'
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '9.67'
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).NumberFormat = '@'
$ws.Cells.Item(27, 5).Value = '  +5.92%  '
$ws.Cells.Item(27, 5).Style = 'Normal'
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '160.64'
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).NumberFormat = '@'
$ws.Cells.Item(28, 5).Value = '  -1.65%  '
$ws.Cells.Item(28, 5).Style = 'Normal'
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '20.11'
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Cells.Item(29, 5).NumberFormat = '@'
$ws.Cells.Item(29, 5).Value = '  +3.67%  '
$ws.Cells.Item(29, 5).Style = 'Normal'
$ws.Cells.Item(30, 5).NumberFormat = '@'
$ws.Cells.Item(30, 5).Value = '  +29.22%  '
$ws.Cells.Item(30, 5).Style = 'Normal'
$ws.Cells.Item(31, 5).NumberFormat = '@'
$ws.Cells.Item(31, 5).Value = '  +2.83%  '
$ws.Cells.Item(31, 5).Style = 'Normal'
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '5.22'
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).NumberFormat = '@'
$ws.Cells.Item(32, 5).Value = '  +9.05%  '
$ws.Cells.Item(32, 5).Style = 'Normal'
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '1.22'
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).NumberFormat = '@'
$ws.Cells.Item(33, 5).Value = '  +7.79%  '
$ws.Cells.Item(33, 5).Style = 'Normal'
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '4.73'
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).NumberFormat = '@'
$ws.Cells.Item(34, 5).Value = '  +10.82%  '
$ws.Cells.Item(34, 5).Style = 'Normal'
$ws.Cells.Item(35, 5).NumberFormat = '@'
$ws.Cells.Item(35, 5).Value = '  +5.23%  '
$ws.Cells.Item(35, 5).Style = 'Normal'
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '2.45'
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).NumberFormat = '@'
$ws.Cells.Item(36, 5).Value = '  +3.23%  '
$ws.Cells.Item(36, 5).Style = 'Normal'
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '1.88'
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).NumberFormat = '@'
$ws.Cells.Item(37, 5).Value = '  +4.33%  '
$ws.Cells.Item(37, 5).Style = 'Normal'
$ws.Cells.Item(38, 5).NumberFormat = '@'
$ws.Cells.Item(38, 5).Value = '  -0.10%  '
$ws.Cells.Item(38, 5).Style = 'Normal'
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '6.11'
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).NumberFormat = '@'
$ws.Cells.Item(39, 5).Value = '  +24.00%  '
$ws.Cells.Item(39, 5).Style = 'Normal'
$ws.Cells.Item(40, 5).NumberFormat = '@'
$ws.Cells.Item(40, 5).Value = '  +17.44%  '
$ws.Cells.Item(40, 5).Style = 'Normal'
$ws.Cells.Item(41, 5).NumberFormat = '@'
$ws.Cells.Item(41, 5).Value = '  +25.79%  '
$ws.Cells.Item(41, 5).Style = 'Normal'
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '1.23'
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).NumberFormat = '@'
$ws.Cells.Item(42, 5).Value = '  +4.31%  '
$ws.Cells.Item(42, 5).Style = 'Normal'
$ws.Cells.Item(43, 5).NumberFormat = '@'
$ws.Cells.Item(43, 5).Value = '  +5.09%  '
$ws.Cells.Item(43, 5).Style = 'Normal'
$ws.Cells.Item(44, 2).NumberFormat = '@'
$ws.Cells.Item(44, 2).Value = 'ARBITRUM'
$ws.Cells.Item(44, 2).Style = 'Normal'
$ws.Cells.Item(44, 3).NumberFormat = '@'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(44, 3).Style = 'Normal'
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '1.14'
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).NumberFormat = '@'
$ws.Cells.Item(44, 5).Value = '  +5.91%  '
$ws.Cells.Item(44, 5).Style = 'Normal'
$ws.Cells.Item(45, 2).NumberFormat = '@'
$ws.Cells.Item(45, 2).Value = 'HuobiToken'
$ws.Cells.Item(45, 2).Style = 'Normal'
$ws.Cells.Item(45, 3).NumberFormat = '@'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(45, 3).Style = 'Normal'
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '2.91'
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).NumberFormat = '@'
$ws.Cells.Item(45, 5).Value = '  +2.85%  '
$ws.Cells.Item(45, 5).Style = 'Normal'
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '17.02'
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).NumberFormat = '@'
$ws.Cells.Item(46, 5).Value = '  +9.61%  '
$ws.Cells.Item(46, 5).Style = 'Normal'
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '7.99'
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).NumberFormat = '@'
$ws.Cells.Item(47, 5).Value = '  +6.81%  '
$ws.Cells.Item(47, 5).Style = 'Normal'
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '95.53'
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).NumberFormat = '@'
$ws.Cells.Item(48, 5).Value = '  +5.15%  '
$ws.Cells.Item(48, 5).Style = 'Normal'
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '1.416.96'
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).NumberFormat = '@'
$ws.Cells.Item(49, 5).Value = '  +3.12%  '
$ws.Cells.Item(49, 5).Style = 'Normal'
$ws.Cells.Item(50, 5).NumberFormat = '@'
$ws.Cells.Item(50, 5).Value = '  +2.56%  '
$ws.Cells.Item(50, 5).Style = 'Normal'
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '47.42'
$ws.Cells.Item(51, 4).Style = 'Normal'
$ws.Cells.Item(51, 5).NumberFormat = '@'
$ws.Cells.Item(51, 5).Value = '  +4.01%  '
$ws.Cells.Item(51, 5).Style = 'Normal'
